$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for columns G:T (Ligand/Receptor expression & derived
# specificity metrics), rows 2-26. Each entry: row number, then the 14 values for
# G..T in order.
$rowData = @(
    @(2, 5.404689, 16.214067, 0.08747859311663772, 0.09021076583983562, 3, 1, 0.6435283333333334, 1.930585, 0.1384187503011309, 0.1430362385488992, 3.478070504355, 31.302634539195, 0.0121086775373061, 0.01290340862234561),
    @(3, 5.404689, 16.214067, 0.08747859311663772, 0.09021076583983562, 3, 1, 1.080588666666667, 3.241766, 0.2324275794584003, 0.2401810927235583, 5.840245680258001, 52.562211122322, 0.02033243765252639, 0.02166692031484076),
    @(4, 5.404689, 16.214067, 0.08747859311663772, 0.09021076583983562, 3, 1, 1.770781666666667, 5.312345000000001, 0.3808835954223518, 0.3935894284240538, 9.570524195235002, 86.13471775711501, 0.03331916106875397, 0.03550600376459706),
    @(5, 5.404689, 16.214067, 0.08747859311663772, 0.09021076583983562, 2, 1, 0.45025, 0.9005000000000001, 0.09684584049355638, 0.06671766993594362, 2.43346122225, 14.6007673335, 0.008471937875574615, 0.00601865209997085),
    @(6, 5.404689, 16.214067, 0.08747859311663772, 0.09021076583983562, 3, 1, 0.7039926666666667, 2.111978, 0.1514242343245606, 0.1564755703675451, 3.804861421614, 34.243752794526, 0.01324637898247664, 0.01411578103808133),
    @(7, 15.12977766666667, 45.389333, 0.2448858138641327, 0.2525342032254661, 3, 1, 0.6435283333333334, 1.930585, 0.1384187503011309, 0.1430362385488992, 9.73644060553389, 87.627965449805, 0.0338967883215486, 0.03612154253431395),
    @(8, 15.12977766666667, 45.389333, 0.2448858138641327, 0.2525342032254661, 3, 1, 1.080588666666667, 3.241766, 0.2324275794584003, 0.2401810927235583, 16.34906627578644, 147.141596482078, 0.05691821696014075, 0.06065394088076558),
    @(9, 15.12977766666667, 45.389333, 0.2448858138641327, 0.2525342032254661, 3, 1, 1.770781666666667, 5.312345000000001, 0.3808835954223518, 0.3935894284240538, 26.79153291287612, 241.123796215885, 0.0932729892524997, 0.09939479270503505),
    @(10, 15.12977766666667, 45.389333, 0.2448858138641327, 0.2525342032254661, 2, 1, 0.45025, 0.9005000000000001, 0.09684584049355638, 0.06671766993594362, 6.812182394416667, 40.87309436650001, 0.02371617246862054, 0.01684849361833316),
    @(11, 15.12977766666667, 45.389333, 0.2448858138641327, 0.2525342032254661, 3, 1, 0.7039926666666667, 2.111978, 0.1514242343245606, 0.1564755703675451, 10.65125252563045, 95.861272730674, 0.03708164686132317, 0.03951543348701835),
    @(12, 17.564497, 52.69349099999999, 0.2842934138044583, 0.2931725999334087, 3, 1, 0.6435283333333334, 1.930585, 0.1384187503011309, 0.1430362385488992, 11.30325148024833, 101.729263322235, 0.03935153905765539, 0.04193430594007604),
    @(13, 17.564497, 52.69349099999999, 0.2842934138044583, 0.2931725999334087, 3, 1, 1.080588666666667, 3.241766, 0.2324275794584003, 0.2401810927235583, 18.97999639390067, 170.819967545106, 0.06607763002653562, 0.0704145154086127),
    @(14, 17.564497, 52.69349099999999, 0.2842934138044583, 0.2931725999334087, 3, 1, 1.770781666666667, 5.312345000000001, 0.3808835954223518, 0.3935894284240538, 31.10288927182167, 279.926003446395, 0.1082826976047365, 0.1153896360373841),
    @(15, 17.564497, 52.69349099999999, 0.2842934138044583, 0.2931725999334087, 2, 1, 0.45025, 0.9005000000000001, 0.09684584049355638, 0.06671766993594362, 7.908414774250001, 47.4504886455, 0.02753263460667519, 0.01955979275661961),
    @(16, 17.564497, 52.69349099999999, 0.2842934138044583, 0.2931725999334087, 3, 1, 0.7039926666666667, 2.111978, 0.1514242343245606, 0.1564755703675451, 12.36527708168867, 111.287493735198, 0.04304891250885556, 0.04587434979071624),
    @(17, 5.613580499999999, 11.227161, 0.09085964511315853, 0.06246494429911598, 3, 1, 0.6435283333333334, 1.930585, 0.1384187503011309, 0.1430362385488992, 3.6124981031975, 21.674988619185, 0.01257667852936766, 0.008934750673712055),
    @(18, 5.613580499999999, 11.227161, 0.09085964511315853, 0.06246494429911598, 3, 1, 1.080588666666667, 3.241766, 0.2324275794584003, 0.2401810927235583, 6.065971467720999, 36.395828806326, 0.02111828738410071, 0.01500289857867788),
    @(19, 5.613580499999999, 11.227161, 0.09085964511315853, 0.06246494429911598, 3, 1, 1.770781666666667, 5.312345000000001, 0.3808835954223518, 0.3935894284240538, 9.9404254337575, 59.642552602545, 0.03460694830949874, 0.02458554172322942),
    @(20, 5.613580499999999, 11.227161, 0.09085964511315853, 0.06246494429911598, 2, 1, 0.45025, 0.9005000000000001, 0.09684584049355638, 0.06671766993594362, 2.527514620125, 10.1100584805, 0.008799378697930091, 0.004167515536315523),
    @(21, 5.613580499999999, 11.227161, 0.09085964511315853, 0.06246494429911598, 3, 1, 0.7039926666666667, 2.111978, 0.1514242343245606, 0.1564755703675451, 3.951919505743, 23.711517034458, 0.01375835219226133, 0.009774237787181107),
    @(22, 18.07044533333334, 54.211336, 0.2924825341016128, 0.3016174867021735, 3, 1, 0.6435283333333334, 1.930585, 0.1384187503011309, 0.1430362385488992, 11.62884356795111, 104.65959211156, 0.04048506685525314, 0.04314223077845152),
    @(23, 18.07044533333334, 54.211336, 0.2924825341016128, 0.3016174867021735, 3, 1, 1.080588666666667, 3.241766, 0.2324275794584003, 0.2401810927235583, 19.52671842881956, 175.740465859376, 0.0679810074350969, 0.07244281754066134),
    @(24, 18.07044533333334, 54.211336, 0.2924825341016128, 0.3016174867021735, 3, 1, 1.770781666666667, 5.312345000000001, 0.3808835954223518, 0.3935894284240538, 31.9988133047689, 287.98931974292, 0.1114017991868629, 0.1187134541938081),
    @(25, 18.07044533333334, 54.211336, 0.2924825341016128, 0.3016174867021735, 2, 1, 0.45025, 0.9005000000000001, 0.09684584049355638, 0.06671766993594362, 8.136218011333336, 48.81730806800001, 0.02832571684475596, 0.02012321592470448),
    @(26, 18.07044533333334, 54.211336, 0.2924825341016128, 0.3016174867021735, 3, 1, 0.7039926666666667, 2.111978, 0.1514242343245606, 0.1564755703675451, 12.72146099806756, 114.493148982608, 0.0442889437796439, 0.04719576826454805)
)

foreach ($entry in $rowData) {
    $row = $entry[0]
    $values = $entry[1..14]
    $arr = New-Object "object[,]" 1,14
    for ($i = 0; $i -lt 14; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range("G$row`:T$row").Value2 = $arr
}
